$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1) Remove the two rows that were dropped entirely from the table:
#    "RM 232" (row 26) and "SC 92" (row 28, becomes row 27 once the first
#    row above it is removed). Deleting top-down keeps the row numbers
#    simple since row 26 is above row 28.
# ----------------------------------------------------------------------
$ws.Rows(26).Delete()   # was "RM 232"; "SC 92" shifts from 28 -> 27
$ws.Rows(27).Delete()   # was "SC 92"

# ----------------------------------------------------------------------
# 2) Cell-level value changes within the rows that stayed in place
#    (rows 1-25 are untouched by the row deletions above).
# ----------------------------------------------------------------------
$ws.Range("C3").Value = ""        # 11.2  -> blank (now missing)
$ws.Range("F5").Value = ""        # 17.66 -> blank (now missing)
$ws.Range("F8").Value = 17.05     # blank -> 17.05
$ws.Range("F10").Value = 16.43    # blank -> 16.43
$ws.Range("F12").Value = ""       # 17.45 -> blank (now missing)
$ws.Range("F15").Value = 16.2     # blank -> 16.2
$ws.Range("F18").Value = ""       # 18.35 -> blank (now missing)
$ws.Range("F19").Value = ""       # 17.81 -> blank (now missing)
$ws.Range("F25").Value = 16.6     # blank -> 16.6

# ----------------------------------------------------------------------
# 3) Cell-level value changes within the rows that shifted up after the
#    two deletions above (final row numbers, post-delete):
#      row 26 -> "SC 5"
#      row 27 -> "SC 101"
#      row 29 -> "SC 119"
#      row 33 -> "SC 232"
# ----------------------------------------------------------------------
$ws.Range("B26").Value = -20.2    # SC 5:   blank -> -20.2
$ws.Range("B27").Value = ""       # SC 101: -20.4 -> blank (now missing)
$ws.Range("F29").Value = ""       # SC 119: 18.06 -> blank (now missing)
$ws.Range("B33").Value = -19.5    # SC 232: blank -> -19.5
$ws.Range("C33").Value = 10.4     # SC 232: blank -> 10.4
